# A new contact (+5521985109311 / DDD 21 / inscribed 2024-09-09) is added
# at the top of the list. Insert a fresh row 7 and push the existing
# rows 7-10 down to 8-11 (dimension grows from A1:C10 to A1:C11).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Insert()

# The phone number, DDD and date columns are stored as plain text in this
# sheet (e.g. "+553291004823"), so writing the new values straight into
# Range.Value would let Excel's smart-parsing turn "+5521985109311" into a
# number and "2024-09-09" into a date serial. Stage the literal text in a
# scratch area with an explicit text format, then copy/paste-values into
# the target cells so the destination keeps its original (General) style
# while still receiving plain text content.
$scratch = $ws.Range("Z1:Z3")
$scratch.NumberFormat = "@"
$ws.Range("Z1").Value = "+5521985109311"
$ws.Range("Z2").Value = "21"
$ws.Range("Z3").Value = "2024-09-09"

$ws.Range("Z1").Copy()
$ws.Range("A7").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("Z2").Copy()
$ws.Range("B7").PasteSpecial(-4163)

$ws.Range("Z3").Copy()
$ws.Range("C7").PasteSpecial(-4163)

$scratch.Clear()
